$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E1 text change ("Stellgröße xN in %" -> "Stellgröße xN") ---
$ws.Range("E1").Value = "Stellgröße xN"

# --- Row 2 ---
$ws.Range("E2").Value = 0.2

# --- Row 3 ---
$ws.Range("B3").Value = 900
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0

# --- Row 4 ---
$ws.Range("C4").Value = 0.2
$ws.Range("E4").Value = 0.2

# --- Row 5 ---
$ws.Range("E5").Value = 0

# --- Row 6 ---
$ws.Range("E6").Value = -0.6

# --- Row 7 ---
$ws.Range("C7").Value = 0.7
$ws.Range("D7").Value = 0.8
$ws.Range("E7").Value = 0.7

# --- Row 8 ---
$ws.Range("E8").Value = -1

# --- Row 9 ---
$ws.Range("E9").Value = 1

# --- New row 10 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 700
$ws.Range("C10").Value = -1
$ws.Range("D10").Value = -1
$ws.Range("E10").Value = -1
$ws.Range("F10").Value = 1
$ws.Range("G10").Formula = "=4*B10"
$ws.Range("H10").Formula = "=5*E10"
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0

# --- Row 11 becomes full data row (replacing the near-empty formatted row) ---
# Clear the old formatting (green Courier style) from D11:F11 first
$ws.Range("D11:F11").ClearFormats()
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 25000
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Formula = "=4*B11"
$ws.Range("H11").Formula = "=5*E11"
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0

# --- Selection ---
$ws.Range("E3").Select() | Out-Null
